$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Sin SmartScore"
$ws.Range("B4").Value = "Con SmartScore"
$ws.Range("B7").Value = "Sin SmartScore"
$ws.Range("B9").Value = "Con SmartScore"
$ws.Range("B14").Value = "Con SmartScore"
$ws.Range("B17").Value = "Sin SmartScore"
$ws.Range("B18").Value = "Sin SmartScore"
$ws.Range("B20").Value = "Con SmartScore"
$ws.Range("B21").Value = "Con SmartScore"
$ws.Range("B22").Value = "Con SmartScore"
$ws.Range("B23").Value = "Con SmartScore"
$ws.Range("B24").Value = "Sin SmartScore"

$ws.Range("I24").Value = 0.56
$ws.Range("L24").Value = 0.463
$ws.Range("O24").Value = 0.43
$ws.Range("R24").Value = 0.718
$ws.Range("U24").Value = 0.584
$ws.Range("X24").Value = 0.561
$ws.Range("AD24").Value = 0.517
$ws.Range("AG24").Value = 0.504
